$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = '0.63 / 0.69'
$ws.Range("C2").Value = '0.72 / 0.83'
$ws.Range("D2").Value = '0.42 / 0.49'
$ws.Range("B3").Value = '0.47 / 0.57'
$ws.Range("C3").Value = '0.46 / 0.58'
$ws.Range("D3").Value = '0.38 / 0.46'
$ws.Range("B4").Value = '0.60 (0.58 - 0.62) / 0.67 (0.65 - 0.68)'
$ws.Range("C4").Value = '0.73 (0.72 - 0.74) / 0.83 (0.83 - 0.84)'
$ws.Range("D4").Value = '0.40 (0.37 - 0.45) / 0.47 (0.44 - 0.50)'
$ws.Range("B5").Value = '0.46 (0.44 - 0.48) / 0.55 (0.52 - 0.58)'
$ws.Range("C5").Value = '0.46 (0.44 - 0.50) / 0.58 (0.56 - 0.61)'
$ws.Range("D5").Value = '0.40 (0.37 - 0.43) / 0.47 (0.45 - 0.49)'
$ws.Range("B6").Value = '0.58 / 0.67'
$ws.Range("C6").Value = '0.75 / 0.85'
$ws.Range("D6").Value = '0.39 / 0.53'
$ws.Range("B7").Value = '0.51 / 0.62'
$ws.Range("C7").Value = '0.65 / 0.75'
$ws.Range("D7").Value = '0.39 / 0.50'
$ws.Range("B8").Value = '0.58 (0.58 - 0.58) / 0.66 (0.66 - 0.67)'
$ws.Range("C8").Value = '0.74 (0.74 - 0.75) / 0.84 (0.82 - 0.85)'
$ws.Range("D8").Value = '0.36 (0.32 - 0.39) / 0.47 (0.43 - 0.52)'
$ws.Range("B9").Value = '0.57 (0.53 - 0.59) / 0.67 (0.65 - 0.69)'
$ws.Range("C9").Value = '0.61 (0.59 - 0.64) / 0.70 (0.68 - 0.73)'
$ws.Range("D9").Value = '0.40 (0.40 - 0.41) / 0.51 (0.51 - 0.52)'
$ws.Range("B10").Value = '0.63 / 0.69'
$ws.Range("C10").Value = '0.79 / 0.87'
$ws.Range("D10").Value = '0.43 / 0.51'
$ws.Range("B11").Value = '0.45 / 0.55'
$ws.Range("C11").Value = '0.35 / 0.43'
$ws.Range("D11").Value = '0.31 / 0.41'
$ws.Range("B12").Value = '0.55 (0.48 - 0.63) / 0.61 (0.52 - 0.70)'
$ws.Range("C12").Value = '0.41 (0.23 - 0.75) / 0.50 (0.32 - 0.84)'
$ws.Range("D12").Value = '0.41 (0.31 - 0.47) / 0.50 (0.41 - 0.55)'
$ws.Range("B13").Value = '0.39 (0.32 - 0.42) / 0.48 (0.42 - 0.52)'
$ws.Range("C13").Value = '0.45 (0.41 - 0.50) / 0.54 (0.49 - 0.60)'
$ws.Range("D13").Value = '0.28 (0.26 - 0.31) / 0.38 (0.37 - 0.38)'
$ws.Range("B14").Value = '0.56 / 0.67'
$ws.Range("C14").Value = '0.72 / 0.84'
$ws.Range("D14").Value = '0.39 / 0.49'
$ws.Range("B15").Value = '0.47 / 0.59'
$ws.Range("C15").Value = '0.65 / 0.74'
$ws.Range("D15").Value = '0.47 / 0.55'
$ws.Range("B16").Value = '0.64 (0.59 - 0.67) / 0.71 (0.69 - 0.72)'
$ws.Range("C16").Value = '0.61 (0.35 - 0.77) / 0.72 (0.45 - 0.88)'
$ws.Range("D16").Value = '0.48 (0.45 - 0.50) / 0.55 (0.51 - 0.57)'
$ws.Range("B17").Value = '0.53 (0.53 - 0.54) / 0.67 (0.66 - 0.67)'
$ws.Range("C17").Value = '0.64 (0.60 - 0.67) / 0.72 (0.68 - 0.75)'
$ws.Range("D17").Value = '0.41 (0.40 - 0.44) / 0.52 (0.50 - 0.53)'
$ws.Range("B18").Value = '0.60 / 0.70'
$ws.Range("C18").Value = '0.77 / 0.85'
$ws.Range("D18").Value = '0.46 / 0.55'
$ws.Range("B19").Value = '0.62 / 0.72'
$ws.Range("C19").Value = '0.83 / 0.87'
$ws.Range("D19").Value = '0.45 / 0.56'
$ws.Range("B20").Value = '0.60 (0.57 - 0.64) / 0.69 (0.68 - 0.71)'
$ws.Range("C20").Value = '0.77 (0.76 - 0.78) / 0.86 (0.83 - 0.87)'
$ws.Range("D20").Value = '0.47 (0.41 - 0.52) / 0.54 (0.50 - 0.58)'
$ws.Range("B21").Value = '0.62 (0.59 - 0.64) / 0.72 (0.71 - 0.73)'
$ws.Range("C21").Value = '0.82 (0.81 - 0.83) / 0.87 (0.86 - 0.87)'
$ws.Range("D21").Value = '0.46 (0.44 - 0.47) / 0.56 (0.55 - 0.58)'
$ws.Range("B22").Value = '0.66 / 0.70'
$ws.Range("C22").Value = '0.70 / 0.82'
$ws.Range("D22").Value = '0.35 / 0.42'
$ws.Range("B23").Value = '0.46 / 0.52'
$ws.Range("C23").Value = '0.47 / 0.60'
$ws.Range("D23").Value = '0.21 / 0.30'
$ws.Range("B24").Value = '0.59 (0.55 - 0.62) / 0.65 (0.62 - 0.67)'
$ws.Range("C24").Value = '0.71 (0.69 - 0.73) / 0.82 (0.81 - 0.83)'
$ws.Range("D24").Value = '0.29 (0.28 - 0.31) / 0.35 (0.34 - 0.36)'
$ws.Range("B25").Value = '0.52 (0.49 - 0.56) / 0.57 (0.52 - 0.61)'
$ws.Range("C25").Value = '0.50 (0.43 - 0.54) / 0.61 (0.54 - 0.66)'
$ws.Range("D25").Value = '0.24 (0.23 - 0.25) / 0.31 (0.30 - 0.32)'
$ws.Range("B26").Value = '0.64 / 0.73'
$ws.Range("C26").Value = '0.79 / 0.86'
$ws.Range("D26").Value = '0.43 / 0.53'
$ws.Range("B27").Value = '0.65 / 0.73'
$ws.Range("C27").Value = '0.76 / 0.82'
$ws.Range("D27").Value = '0.59 / 0.65'
$ws.Range("B28").Value = '0.63 (0.60 - 0.65) / 0.73 (0.70 - 0.77)'
$ws.Range("C28").Value = '0.70 (0.52 - 0.83) / 0.76 (0.59 - 0.87)'
$ws.Range("D28").Value = '0.48 (0.47 - 0.48) / 0.57 (0.56 - 0.59)'
$ws.Range("B29").Value = '0.66 (0.65 - 0.67) / 0.76 (0.75 - 0.77)'
$ws.Range("C29").Value = '0.72 (0.70 - 0.76) / 0.78 (0.75 - 0.83)'
$ws.Range("D29").Value = '0.53 (0.49 - 0.59) / 0.60 (0.55 - 0.66)'
$ws.Range("B30").Value = '0.60 / 0.70'
$ws.Range("C30").Value = '0.73 / 0.82'
$ws.Range("D30").Value = '0.38 / 0.48'
$ws.Range("B31").Value = '0.50 / 0.62'
$ws.Range("C31").Value = '0.71 / 0.77'
$ws.Range("D31").Value = '0.31 / 0.43'
$ws.Range("B32").Value = '0.60 (0.59 - 0.61) / 0.71 (0.70 - 0.71)'
$ws.Range("C32").Value = '0.73 (0.72 - 0.75) / 0.82 (0.81 - 0.83)'
$ws.Range("D32").Value = '0.28 (0.25 - 0.34) / 0.40 (0.37 - 0.45)'
$ws.Range("B33").Value = '0.54 (0.49 - 0.57) / 0.67 (0.62 - 0.69)'
$ws.Range("C33").Value = '0.63 (0.61 - 0.67) / 0.71 (0.68 - 0.75)'
$ws.Range("D33").Value = '0.30 (0.29 - 0.32) / 0.42 (0.41 - 0.45)'
$ws.Range("B34").Value = '0.70 / 0.75'
$ws.Range("C34").Value = '0.46 / 0.56'
$ws.Range("D34").Value = '0.38 / 0.47'
$ws.Range("B35").Value = '0.46 / 0.57'
$ws.Range("C35").Value = '0.64 / 0.74'
$ws.Range("D35").Value = '0.45 / 0.55'
$ws.Range("B36").Value = '0.64 (0.60 - 0.71) / 0.71 (0.67 - 0.74)'
$ws.Range("C36").Value = '0.57 (0.16 - 0.78) / 0.65 (0.25 - 0.87)'
$ws.Range("D36").Value = '0.34 (0.27 - 0.42) / 0.44 (0.40 - 0.51)'
$ws.Range("B37").Value = '0.58 (0.52 - 0.64) / 0.68 (0.61 - 0.73)'
$ws.Range("C37").Value = '0.63 (0.57 - 0.70) / 0.74 (0.67 - 0.80)'
$ws.Range("D37").Value = '0.48 (0.46 - 0.49) / 0.58 (0.56 - 0.61)'
$ws.Range("B38").Value = '0.21 / 0.28'
$ws.Range("C38").Value = '0.38 / 0.46'
$ws.Range("D38").Value = '0.31 / 0.35'
$ws.Range("B39").Value = '0.38 / 0.43'
$ws.Range("C39").Value = '0.20 / 0.27'
$ws.Range("B40").Value = '0.43 (0.21 - 0.56) / 0.52 (0.29 - 0.65)'
$ws.Range("C40").Value = '0.15 (0.14 - 0.16) / 0.25 (0.23 - 0.27)'
$ws.Range("D40").Value = '0.32 (0.29 - 0.36) / 0.40 (0.34 - 0.43)'
$ws.Range("B41").Value = '0.47 (0.37 - 0.53) / 0.53 (0.43 - 0.60)'
$ws.Range("C41").Value = '0.43 (0.36 - 0.50) / 0.52 (0.48 - 0.59)'
$ws.Range("D41").Value = '0.36 (0.33 - 0.38) / 0.44 (0.41 - 0.48)'
